$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row reorder: Algorand/BinanceUSD swap (rows 42-43) ---
$ws.Range("B42").Value = "Algorand"
$ws.Range("C42").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.201"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +19.23%  "
$ws.Range("B43").Value = "BinanceUSD"
$ws.Range("C43").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.00"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.01%  "

# --- Row reorder: Cronos/NEARProtocol/ARBITRUM rotation (rows 45-47) ---
$ws.Range("B45").Value = "Cronos"
$ws.Range("C45").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.101"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +5.33%  "
$ws.Range("B46").Value = "NEARProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.51"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +14.23%  "
$ws.Range("B47").Value = "ARBITRUM"
$ws.Range("C47").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.20"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.71%  "

# --- Price / volume updates ---
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.900.99"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.58%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.363.82"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.35%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.689"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +6.11%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "241.20"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.06%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "76.44"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +7.44%  "
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.625"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +27.35%  "
$ws.Range("E10").Value = "  +5.40%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "57.36"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.88%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "32.88"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +20.89%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.52"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +19.87%  "
$ws.Range("E14").Value = "  +1.52%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.719.44"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.39%  "
$ws.Range("E16").Value = "  +5.24%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.924"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +6.93%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.369.77"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.84%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "43.886.50"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.47%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0000102"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.80%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.65"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +5.14%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "77.72"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.58%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "257.20"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.07%  "
$ws.Range("E24").Value = "  -0.03%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.54"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.38%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.17"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +10.83%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.67"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.89%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.76"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +16.93%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.30"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.74%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "23.19"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.80%  "
$ws.Range("E31").Value = "  +1.84%  "
$ws.Range("E32").Value = "  -0.92%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.134"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +6.29%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.33"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +7.13%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0753"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +9.45%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.40"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +6.93%  "
$ws.Range("E37").Value = "  +4.16%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.46"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.46%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.50"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.47%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0274"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +7.75%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "19.03"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.82%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.89"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.13%  "
$ws.Range("E48").Value = "  +4.23%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "102.30"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.46%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "4.50"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.52%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "54.72"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +9.01%  "
